$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.953.23'
$ws.Range("D3").Value = '1.640.19'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '213.04'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").Value = '0.524'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '23.59'
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = '0.0881'
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").Value = '1.872.80'
$ws.Range("D13").Value = '1.637.13'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("D15").Value = '0.574'
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("D16").Value = '65.89'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '27.956.45'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").Value = '233.89'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '10.75'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '151.19'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").Value = '6.99'
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  +1.77%  '
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").Value = '1.411.24'
$ws.Range("E34").Value = '  -4.70%  '
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("D37").Value = '0.884'
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").Value = '0.558'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  -5.42%  '
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").Value = '1.89'
$ws.Range("E43").Value = '  +7.58%  '
$ws.Range("D44").Value = '66.54'
$ws.Range("E44").Value = '  -2.19%  '
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = '1.781.89'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '87.93'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").Value = '7.63'
$ws.Range("E51").Value = '  -1.44%  '
